# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("fuentes-mineromedicinales") was previously described as a
# dimension; it is now curated as a measure.
$ws.Range("C2").Value = "iaest-measure:fuentes-mineromedicinales"
$ws.Range("C3").Value = "medida"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("C5").Clear()

# Column E ("municipio-nombre") was previously described as a measure;
# it is now curated as a dimension (refArea).
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"
